# Updates the cryptos list (Price / Volume(1h) columns) to refreshed values.
# D/E columns hold plain text (not numbers), so numeric-looking prices are
# forced to Text format before assignment to keep them as strings, matching
# the original inlineStr cell type (avoids Excel auto-converting "1.001" etc.
# into a float cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.795.83"
$ws.Range("E2").Value = "  -4.14%  "

$ws.Range("D3").Value = "1.720.96"
$ws.Range("E3").Value = "  -2.64%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue $ws.Range("D5") "308.81"
$ws.Range("E5").Value = "  -6.00%  "

Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.09%  "

Set-TextValue $ws.Range("D7") "0.4889"
$ws.Range("E7").Value = "  +4.91%  "

Set-TextValue $ws.Range("D8") "0.3484"
$ws.Range("E8").Value = "  -1.01%  "

Set-TextValue $ws.Range("D9") "42.73"
$ws.Range("E9").Value = "  -1.53%  "

Set-TextValue $ws.Range("D10") "0.07221"
$ws.Range("E10").Value = "  -2.20%  "

Set-TextValue $ws.Range("D11") "1.046"
$ws.Range("E11").Value = "  -3.52%  "

Set-TextValue $ws.Range("D12") "1.001"
$ws.Range("E12").Value = "  +0.08%  "

Set-TextValue $ws.Range("D13") "19.76"
$ws.Range("E13").Value = "  -4.16%  "

Set-TextValue $ws.Range("D14") "5.852"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "1.720.16"
$ws.Range("E15").Value = "  -2.59%  "

Set-TextValue $ws.Range("D16") "6.797"
$ws.Range("E16").Value = "  -5.41%  "

Set-TextValue $ws.Range("D17") "86.27"
$ws.Range("E17").Value = "  -6.44%  "

Set-TextValue $ws.Range("D18") "0.00001034"
$ws.Range("E18").Value = "  -2.05%  "

Set-TextValue $ws.Range("D19") "0.06404"
$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("E20").Value = "  +0.14%  "

Set-TextValue $ws.Range("D21") "16.48"
$ws.Range("E21").Value = "  -2.58%  "

Set-TextValue $ws.Range("D22") "5.699"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").Value = "26.852.23"
$ws.Range("E23").Value = "  -4.06%  "

Set-TextValue $ws.Range("D24") "10.91"
$ws.Range("E24").Value = "  -1.97%  "

Set-TextValue $ws.Range("D25") "2.052"
$ws.Range("E25").Value = "  -4.80%  "

Set-TextValue $ws.Range("D26") "154.20"
$ws.Range("E26").Value = "  -6.13%  "

Set-TextValue $ws.Range("D27") "19.78"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "1.917.26"
$ws.Range("E28").Value = "  -2.61%  "

Set-TextValue $ws.Range("D29") "2.059"
$ws.Range("E29").Value = "  -6.02%  "

Set-TextValue $ws.Range("D30") "119.72"
$ws.Range("E30").Value = "  -2.74%  "

Set-TextValue $ws.Range("D31") "1.031"
$ws.Range("E31").Value = "  -4.19%  "

Set-TextValue $ws.Range("D32") "0.09317"
$ws.Range("E32").Value = "  +0.09%  "

Set-TextValue $ws.Range("D33") "3.567"
$ws.Range("E33").Value = "  -2.40%  "

Set-TextValue $ws.Range("D34") "5.337"
$ws.Range("E34").Value = "  -3.84%  "

Set-TextValue $ws.Range("D35") "0.05881"
$ws.Range("E35").Value = "  -3.72%  "

Set-TextValue $ws.Range("D36") "0.02174"
$ws.Range("E36").Value = "  -4.17%  "

Set-TextValue $ws.Range("D37") "1.424"
$ws.Range("E37").Value = "  -1.85%  "

Set-TextValue $ws.Range("D38") "10.92"
$ws.Range("E38").Value = "  -6.53%  "

Set-TextValue $ws.Range("D41") "4.729"
$ws.Range("E41").Value = "  -3.80%  "

Set-TextValue $ws.Range("D42") "0.5948"
$ws.Range("E42").Value = "  -3.43%  "

Set-TextValue $ws.Range("D43") "1.116"
$ws.Range("E43").Value = "  -6.21%  "

Set-TextValue $ws.Range("D44") "7.403"
$ws.Range("E44").Value = "  -4.88%  "

Set-TextValue $ws.Range("D45") "12.81"
$ws.Range("E45").Value = "  -2.67%  "

Set-TextValue $ws.Range("D46") "3.570"
$ws.Range("E46").Value = "  -4.65%  "

Set-TextValue $ws.Range("D47") "0.5579"
$ws.Range("E47").Value = "  -3.77%  "

Set-TextValue $ws.Range("D48") "119.09"
$ws.Range("E48").Value = "  -3.79%  "

Set-TextValue $ws.Range("D49") "1.827"
$ws.Range("E49").Value = "  -5.50%  "

Set-TextValue $ws.Range("D50") "0.06622"
$ws.Range("E50").Value = "  -2.87%  "

Set-TextValue $ws.Range("D51") "1.090"
$ws.Range("E51").Value = "  -3.03%  "

# Rows 39/40: Frax and Algorand swapped ranking positions, with updated values
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D39") "0.1983"
$ws.Range("E39").Value = "  -4.01%  "

$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D40") "1.001"
$ws.Range("E40").Value = "  +0.07%  "
